$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to Wolters_2018 (file B3 = "metrics_sim_with_priors.json")
# Update the metric values per the corrected relevance markers.

$ws.Range("C3").Value = 0.8947368421052632
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("H3").Value = 0.8100467289719626
$ws.Range("I3").Value = 0.03605528724416063
$ws.Range("J3").Value = 0.7894736842105263
$ws.Range("K3").Value = 163.6315789473684

$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 38
$ws.Range("T3").Value = 157
$ws.Range("U3").Value = 317
$ws.Range("V3").Value = 4256
$ws.Range("W3").Value = 4251
$ws.Range("X3").Value = 4223
$ws.Range("Y3").Value = 4104
$ws.Range("Z3").Value = 3944

$ws.Range("AF3").Value = 0.998827
$ws.Range("AG3").Value = 0.997653
$ws.Range("AH3").Value = 0.991082
$ws.Range("AI3").Value = 0.963154
$ws.Range("AJ3").Value = 0.925604
